# Applies the "built volun filter, need modify" edit:
#  - Sheet "Volunteer Details": fixes row 11 (ID 10028) and appends new
#    volunteer rows 12-16 (IDs 10029,10030,10031,10032,10034) plus a
#    trailing "CU / SSS / S" marker row 17.
#  - Sheet "Volunteer Hours": appends the matching weekly-hours rows
#    15-20 for the same volunteer IDs (all zero hours / not started).

$wb2 = $excel.ActiveWorkbook
$wsDetails = $wb2.Sheets("Volunteer Details")
$wsHours   = $wb2.Sheets("Volunteer Hours")

# ---------------------------------------------------------------------
# Sheet 1: "Volunteer Details"
# ---------------------------------------------------------------------

# Row 11 already existed (ID 1 placeholder) - turn it into volunteer 10028
# and correct the Preferred/Surname-group columns (G/H) from "XXXXX" to "cccc".
$wsDetails.Range("A11").Value = 10028
$wsDetails.Range("B11").Value = ""
$wsDetails.Range("C11").Value = "LU invite sent"
$wsDetails.Range("D11").Value = "NA - Last years Vol"
$wsDetails.Range("E11").Value = "No"
$wsDetails.Range("F11").Value = 165654
$wsDetails.Range("G11").Value = "cccc"
$wsDetails.Range("H11").Value = "cccc"
$wsDetails.Range("I11").Value = "ccccc"
$wsDetails.Range("J11").Value = "Female"
$wsDetails.Range("K11").Value = 37145
$wsDetails.Range("L11").Value = "XXXXX"
$wsDetails.Range("O11").Value = "Yes"
$wsDetails.Range("P11").Value = "No"

# Rows 12 and 16 onward repeat the same shape with a couple of value
# swaps (F/G/H/I/O) as new volunteers get appended underneath.
$detailRows = @(
    @{ Row = 12; A = 10029; F = 165654;  GHI = "cccc"; O = "Yes" },
    @{ Row = 13; A = 10030; F = 1677764; GHI = "dddd"; O = "No"  },
    @{ Row = 14; A = 10031; F = 1677764; GHI = "dddd"; O = "No"  },
    @{ Row = 15; A = 10032; F = 1677764; GHI = "dddd"; O = "No"  },
    @{ Row = 16; A = 10034; F = 1677764; GHI = "dddd"; O = "No"  }
)

foreach ($r in $detailRows) {
    $row = $r.Row
    $wsDetails.Range("A$row").Value = $r.A
    $wsDetails.Range("B$row").Value = ""
    $wsDetails.Range("C$row").Value = "LU invite sent"
    $wsDetails.Range("D$row").Value = "NA - Last years Vol"
    $wsDetails.Range("E$row").Value = "No"
    $wsDetails.Range("F$row").Value = $r.F
    $wsDetails.Range("G$row").Value = $r.GHI
    $wsDetails.Range("H$row").Value = $r.GHI
    $wsDetails.Range("I$row").Value = $r.GHI
    $wsDetails.Range("J$row").Value = "Female"
    $wsDetails.Range("K$row").Value = 37145
    $wsDetails.Range("L$row").Value = "XXXXX"
    $wsDetails.Range("O$row").Value = $r.O
    $wsDetails.Range("P$row").Value = "No"
}
# Rows 13-16 use "dddd" across G/H/I (set explicitly, matches GHI above).
$wsDetails.Range("I13").Value = "dddd"
$wsDetails.Range("I14").Value = "dddd"
$wsDetails.Range("I15").Value = "dddd"
$wsDetails.Range("I16").Value = "dddd"

# Trailing marker/totals row.
$wsDetails.Range("A17").Value = 11
$wsDetails.Range("C17").Value = "CU"
$wsDetails.Range("D17").Value = "SSS"
$wsDetails.Range("E17").Value = "S"
$wsDetails.Range("F17").Value = "S"
$wsDetails.Range("G17").Value = "S"
$wsDetails.Range("H17").Value = "S"
$wsDetails.Range("I17").Value = "S"
$wsDetails.Range("J17").Value = "S"
$wsDetails.Range("K17").Value = 37145
$wsDetails.Range("L17").Value = "S"

# ---------------------------------------------------------------------
# Sheet 2: "Volunteer Hours"
# ---------------------------------------------------------------------

$hoursRows = @(
    @{ Row = 15; A = 10028; Label = "cccc"; HasE = $true  },
    @{ Row = 16; A = 10029; Label = "cccc"; HasE = $true  },
    @{ Row = 17; A = 10030; Label = "dddd"; HasE = $false },
    @{ Row = 18; A = 10031; Label = "dddd"; HasE = $true  },
    @{ Row = 19; A = 10032; Label = "dddd"; HasE = $true  },
    @{ Row = 20; A = 10034; Label = "dddd"; HasE = $false }
)

foreach ($r in $hoursRows) {
    $row = $r.Row
    $wsHours.Range("A$row").Value = $r.A
    $wsHours.Range("B$row").Value = $r.Label
    $wsHours.Range("C$row").Value = $r.Label
    if ($r.HasE) {
        $wsHours.Range("E$row").Value = 0
    }
    if ($row -lt 20) {
        $wsHours.Range("F$row").Value = 0
        $wsHours.Range("G$row").Value = 0
        $wsHours.Range("H$row").Value = 0
        $wsHours.Range("I$row").Value = 0
    }
}
